$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 holds account 005055865 / G3C / 628.98 - remove it entirely,
# shifting subsequent rows up.
$ws.Rows.Item(20).Delete()
